# Automatic update of files.
# - Column C ("Förändrad") bumps by one day (46059 -> 46060) for every data row (2-14).
# - Rows 7-14 are re-sorted: the records (identified by columns A/B/G which travel
#   together) are reshuffled into a new row order while D/E (Län/Kommun) stay the
#   same for every row in that block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2-6: only the "Förändrad" date (column C) changes ---
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 3).Value = 46060
}

# --- Rows 7-14: data reordered, column C bumped for all ---

# New row order values for columns A (Beteckning), B (Datum), G (Area (ha)).
$newA = @(
    "A 62884-2021",
    "A 25634-2025",
    "A 28266-2025",
    "A 25015-2023",
    "A 19922-2025",
    "A 60024-2025",
    "A 3402-2026",
    "A 14271-2021"
)

$newB = @(
    44504,
    45803.59570601852,
    45818.56381944445,
    45085.6989699074,
    45771.63034722222,
    45992,
    46042.39047453704,
    44278
)

$newG = @(
    0.8,
    6,
    1.9,
    1.8,
    10.1,
    1.1,
    5.5,
    6.7
)

for ($i = 0; $i -lt 8; $i++) {
    $row = 7 + $i
    $ws.Cells.Item($row, 1).Value = $newA[$i]
    $ws.Cells.Item($row, 2).Value = $newB[$i]
    $ws.Cells.Item($row, 3).Value = 46060
    $ws.Cells.Item($row, 7).Value = $newG[$i]
}
